$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "dummy" row under the header (row 1), using the same
# plain Arial 8pt styling already used elsewhere in the template so the
# resulting cell style matches the existing style palette.
$ws.Range("A2").Value = "dummy"
$ws.Range("A2").Font.Name = "arial"
$ws.Range("A2").Font.Size = 8

# Update the sheet's saved selection/active cell.
$ws.Range("A4").Select()
